# Apply the "hole_id" index column addition to the train sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

# Header for new index column (match style of existing header cells, e.g. B1)
$ws.Range("A1").Value = "hole_id"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# hole_id values for rows 2-29 (in row order)
$holeIds = @(
    "BRG_05_05",
    "ECO_09_02",
    "BRG_01_06",
    "ECO_09_05",
    "ECO_09_04",
    "BRG_01_02",
    "BRG_16_08",
    "BRG_05_11",
    "BRG_01_03",
    "BRG_05_09",
    "BRG_01_08",
    "BRG_05_04",
    "BRG_05_15",
    "ECO_09_01",
    "BRG_01_07",
    "BRG_16_09",
    "BRG_16_05",
    "BRG_05_01",
    "BRG_16_02",
    "BRG_05_03",
    "BRG_05_02",
    "BRG_05_14",
    "BRG_16_07",
    "BRG_08_01",
    "BRG_01_01",
    "BRG_01_09",
    "BRG_01_04",
    "BRG_16_01"
)

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
